# Insert a new row of data for LOC_2018 above the existing LOC_2019 row
# (row 6), shifting LOC_2019..LOC_2023 down by one row each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6:10 down to 7:11 by inserting a new blank row at position 6.
$ws.Rows.Item(6).Insert()

# Force the new row's cells to be treated as plain text, matching the
# inline-string representation used throughout the rest of the sheet,
# so values like dates and numbers aren't auto-converted.
$ws.Range("A6:J6").NumberFormat = "@"

$ws.Range("A6").Value = "LOC_2018"
$ws.Range("B6").Value = "2018-05-01"
$ws.Range("C6").Value = "2018-05-07"
$ws.Range("D6").Value = "249.24"
$ws.Range("E6").Value = "274.43"
$ws.Range("F6").Value = "0.994623944470636"
$ws.Range("G6").Value = "7.20898114739464e-07"
$ws.Range("H6").Value = "4.62177579364686e-05"
$ws.Range("I6").Value = "-70239.7286160114"
$ws.Range("J6").Value = "full_ice_to_functional_ice_off"
